$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-03 Sunday" "2025-08-04 Monday"

Replace-Text "26×15=" "76×65="
Replace-Text "91×27=" "34×38="
Replace-Text "39×29=" "62×88="
Replace-Text "25×33=" "89×38="
Replace-Text "15×51=" "89×52="
Replace-Text "32×97=" "56×49="
Replace-Text "48×82=" "66×19="
Replace-Text "18×74=" "22×82="
Replace-Text "31×58=" "92×26="
Replace-Text "79×26=" "35×43="
Replace-Text "74×49=" "97×44="
Replace-Text "74×34=" "36×43="
Replace-Text "23×31=" "67×45="
Replace-Text "69×34=" "39×46="
Replace-Text "40×56=" "15×68="
Replace-Text "72×63=" "72×94="
Replace-Text "15×97=" "46×66="
Replace-Text "97×24=" "75×42="
Replace-Text "72×23=" "29×29="
Replace-Text "13×91=" "11×18="
Replace-Text "86×30=" "51×38="
Replace-Text "92×76=" "68×14="
Replace-Text "55×50=" "75×72="
Replace-Text "25×71=" "50×85="
Replace-Text "95×49=" "81×40="
